$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 0.76363224855622791
$ws.Range("BD1").Value = 0.90236370739106486
$ws.Range("BP1").Value = 0.78641254315592191
$ws.Range("C2").Value = 0.89766866833559278
$ws.Range("D2").Value = 0.88449701315432716
$ws.Range("AP2").Value = 0.53394532504140146
$ws.Range("BP2").Value = 0.86307471023009574
$ws.Range("E3").Value = 0.78754363863456711
$ws.Range("D5").Value = 0.88758825178109935
$ws.Range("G5").Value = 0.735644638449515
$ws.Range("E6").Value = 0.66774229306390265
$ws.Range("L6").Value = 0.78234391340596943
$ws.Range("F7").Value = 0.92446808820452941
$ws.Range("H7").Value = 0.79211837276729069
$ws.Range("BC7").Value = 0.94598292644579329
$ws.Range("I8").Value = 0.8062211302037352
$ws.Range("J8").Value = 0.94512414790492305
$ws.Range("I10").Value = 0.84979502333526302
$ws.Range("L10").Value = 0.75331026756985731
$ws.Range("I11").Value = 0.86026516273515785
$ws.Range("M11").Value = 0.74113044935761418
$ws.Range("M12").Value = 0.97729665082634021
$ws.Range("N12").Value = 0.58595569507030087
$ws.Range("O13").Value = 0.94175264284272875
$ws.Range("P14").Value = 0.92580098658602283
$ws.Range("N15").Value = 0.801809436362521
$ws.Range("AY15").Value = 0.80089087350759125
$ws.Range("O16").Value = 0.73027181666322827
$ws.Range("BG16").Value = 0.81138439093242054
$ws.Range("O17").Value = 0.91570224681324142
$ws.Range("P17").Value = 0.6958278102736235
$ws.Range("R17").Value = 0.86378532109860884
$ws.Range("P18").Value = 0.90379732261465728
$ws.Range("R19").Value = 0.89177011566850184
$ws.Range("T19").Value = 0.91414688419706225
$ws.Range("AN20").Value = 0.79041807832757538
$ws.Range("J21").Value = 0.72771512266651206
$ws.Range("T21").Value = 0.94117205742077248
$ws.Range("V21").Value = 0.64920798125052781
$ws.Range("R22").Value = 0.95914530182824664
$ws.Range("X22").Value = 0.77053380608181521
$ws.Range("M23").Value = 0.72916232025255179
$ws.Range("V23").Value = 0.88102912955174117
$ws.Range("Y23").Value = 0.63446035126518763
$ws.Range("W24").Value = 0.73296085485563145
$ws.Range("Y24").Value = 0.93022719242919361
$ws.Range("Z24").Value = 0.8746396557331102
$ws.Range("Y26").Value = 0.83010571852257087
$ws.Range("AA26").Value = 0.60952270628391525
$ws.Range("AB26").Value = 0.74852867154837321
$ws.Range("Y27").Value = 0.98571448752947877
$ws.Range("X28").Value = 0.53301407781434018
$ws.Range("AA28").Value = 0.58185811052203074
$ws.Range("AC28").Value = 0.82199447525436353
$ws.Range("AS28").Value = 0.96200493906435813
$ws.Range("AR29").Value = 0.85071600055615315
$ws.Range("AC30").Value = 0.84909768246319195
$ws.Range("AE30").Value = 0.84577422780550959
$ws.Range("AF30").Value = 0.90636583883146615
$ws.Range("AC31").Value = 0.91769353232212247
$ws.Range("AG31").Value = 0.88361427958661065
$ws.Range("AE32").Value = 0.75809484356522705
$ws.Range("AG32").Value = 0.64574147416869399
$ws.Range("AI33").Value = 0.99484501670287317
$ws.Range("AY33").Value = 0.98771998903647762
$ws.Range("AJ34").Value = 0.71236058135964031
$ws.Range("AO34").Value = 0.99983237251179946
$ws.Range("AQ34").Value = 0.75329324056739411
$ws.Range("AA35").Value = 0.99632670357067687
$ws.Range("BE35").Value = 0.82686745320877253
$ws.Range("AI36").Value = 0.99188978246836312
$ws.Range("AK36").Value = 0.80164195442970432
$ws.Range("G37").Value = 0.97148899112603049
$ws.Range("AL37").Value = 0.83335254775500056
$ws.Range("AM38").Value = 0.73011201709607554
$ws.Range("P39").Value = 0.92456767637680481
$ws.Range("AO39").Value = 0.98391317733832395
$ws.Range("AP40").Value = 0.58942741221892891
$ws.Range("AX40").Value = 0.99796374359961937
$ws.Range("AP41").Value = 0.55785539447753218
$ws.Range("F42").Value = 0.74778553698514805
$ws.Range("AR43").Value = 0.94567799777902728
$ws.Range("AS43").Value = 0.86705162951532677
$ws.Range("BN43").Value = 0.94421898466505838
$ws.Range("K45").Value = 0.96930150485740518
$ws.Range("AR45").Value = 0.78868944236625915
$ws.Range("J46").Value = 0.75338274095552027
$ws.Range("AJ46").Value = 0.91376194071068717
$ws.Range("AU46").Value = 0.89040314360114359
$ws.Range("BJ46").Value = 0.65638650849846059
$ws.Range("Q47").Value = 0.79682053387465823
$ws.Range("AS47").Value = 0.96966289071657141
$ws.Range("AV47").Value = 0.89348731287718475
$ws.Range("BA47").Value = 0.63396842437074041
$ws.Range("AT48").Value = 0.87279782971625874
$ws.Range("AW48").Value = 0.91384307524094366
$ws.Range("C49").Value = 0.9177251537139135
$ws.Range("AX49").Value = 0.98193773131751039
$ws.Range("AV50").Value = 0.8166031640016147
$ws.Range("AW51").Value = 0.62784081974336248
$ws.Range("AX51").Value = 0.73328911020265419
$ws.Range("BA51").Value = 0.99566289293985266
$ws.Range("AC52").Value = 0.72556297774852374
$ws.Range("C53").Value = 0.83742916617441088
$ws.Range("BB53").Value = 0.81707105368186617
$ws.Range("AZ54").Value = 0.93789844195905592
$ws.Range("BD54").Value = 0.90903688297191609
$ws.Range("BB55").Value = 0.62523779880553687
$ws.Range("BD55").Value = 0.94911474944831875
$ws.Range("BF56").Value = 0.9475062345021279
$ws.Range("BC57").Value = 0.88698700658657215
$ws.Range("BD57").Value = 0.74440771515520376
$ws.Range("BF57").Value = 0.84969437385096791
$ws.Range("BG57").Value = 0.8739071920124033
$ws.Range("BH58").Value = 0.99268754764503442
$ws.Range("BF59").Value = 0.78923894676607631
$ws.Range("BH59").Value = 0.92105688329855129
$ws.Range("BJ60").Value = 0.90808578142998508
$ws.Range("AT61").Value = 0.88155511486712401
$ws.Range("BG61").Value = 0.91175936157987736
$ws.Range("BH61").Value = 0.61145432910558384
$ws.Range("BJ61").Value = 0.8956381669509641
$ws.Range("BK62").Value = 0.67370187292466033
$ws.Range("AL63").Value = 0.61871614490214033
$ws.Range("AX63").Value = 0.88546282285697486
$ws.Range("BL63").Value = 0.99534393427298073
$ws.Range("AF64").Value = 0.9868701619972815
$ws.Range("Z65").Value = 0.56036572972678189
$ws.Range("AK65").Value = 0.96339768110314239
$ws.Range("BK65").Value = 0.64798238570896083
$ws.Range("BL65").Value = 0.73549426341082214
$ws.Range("B66").Value = 0.68129698890707324
$ws.Range("BL66").Value = 0.71837160985151394
$ws.Range("K67").Value = 0.62346525949232201
$ws.Range("BN67").Value = 0.65126604095195173
$ws.Range("BO68").Value = 0.92118949401695238
